# chore: update Sheets via scheduled runner
# Refreshes market-price-derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leves across the job sheets, as produced by
# the scheduled price-refresh job.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2297.7385
$ws.Range("I15").Value = 2297.7385
$ws.Range("K15").Value = 6893.2155
$ws.Range("M15").Value = -6724.2155
$ws.Range("H41").Value = 477.875
$ws.Range("I41").Value = 373
$ws.Range("J41").Value = 708.6
$ws.Range("K41").Value = 373
$ws.Range("L41").Value = 708.6
$ws.Range("M41").Value = 67
$ws.Range("N41").Value = -1588.6
$ws.Range("H74").Value = 4019.6
$ws.Range("I74").Value = 4019.6
$ws.Range("K74").Value = 4019.6
$ws.Range("M74").Value = -3083.6
$ws.Range("H77").Value = 4019.6
$ws.Range("I77").Value = 4019.6
$ws.Range("K77").Value = 20098
$ws.Range("M77").Value = -15418
$ws.Range("H98").Value = 1910.6296
$ws.Range("I98").Value = 1551.4762
$ws.Range("K98").Value = 1551.4762
$ws.Range("M98").Value = -53.47620000000006
$ws.Range("H108").Value = 41500
$ws.Range("J108").Value = 41500
$ws.Range("L108").Value = 41500
$ws.Range("N108").Value = -49180
$ws.Range("H122").Value = 1910.6296
$ws.Range("I122").Value = 1551.4762
$ws.Range("K122").Value = 4654.4286
$ws.Range("M122").Value = -2204.4286
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 1161.52
$ws.Range("I132").Value = 1183.5454
$ws.Range("K132").Value = 3550.6362
$ws.Range("M132").Value = -1020.6362
$ws.Range("H138").Value = 3418.0356
$ws.Range("I138").Value = 4759.5386
$ws.Range("J138").Value = 2255.4
$ws.Range("K138").Value = 14278.6158
$ws.Range("L138").Value = 6766.200000000001
$ws.Range("M138").Value = -9138.6158
$ws.Range("N138").Value = -17046.2
$ws.Range("H141").Value = 1122626.8
$ws.Range("I141").Value = 1402328.2
$ws.Range("J141").Value = 3820.8
$ws.Range("K141").Value = 4206984.6
$ws.Range("L141").Value = 11462.4
$ws.Range("M141").Value = -4201804.6
$ws.Range("N141").Value = -21822.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2757.8552
$ws.Range("I32").Value = 2125.3936
$ws.Range("J32").Value = 5329.8667
$ws.Range("K32").Value = 2125.3936
$ws.Range("L32").Value = 5329.8667
$ws.Range("M32").Value = -1838.3936
$ws.Range("N32").Value = -5903.8667
$ws.Range("H45").Value = 2901
$ws.Range("I45").Value = 2822.4
$ws.Range("K45").Value = 2822.4
$ws.Range("M45").Value = -2445.4
$ws.Range("H61").Value = 2867
$ws.Range("I61").Value = 1049.3334
$ws.Range("J61").Value = 4684.6665
$ws.Range("K61").Value = 1049.3334
$ws.Range("L61").Value = 4684.6665
$ws.Range("M61").Value = -837.3334
$ws.Range("N61").Value = -5108.6665
$ws.Range("H136").Value = 2867
$ws.Range("I136").Value = 1049.3334
$ws.Range("J136").Value = 4684.6665
$ws.Range("K136").Value = 3148.0002
$ws.Range("L136").Value = 14053.9995
$ws.Range("M136").Value = -598.0001999999999
$ws.Range("N136").Value = -19153.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 47694.25
$ws.Range("J131").Value = 47694.25
$ws.Range("L131").Value = 47694.25
$ws.Range("N131").Value = -57774.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 546.7917
$ws.Range("I107").Value = 459.14285
$ws.Range("J107").Value = 669.5
$ws.Range("K107").Value = 459.14285
$ws.Range("L107").Value = 669.5
$ws.Range("M107").Value = 1460.85715
$ws.Range("N107").Value = -4509.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4495
$ws.Range("J80").Value = 3990
$ws.Range("L80").Value = 11970
$ws.Range("N80").Value = -13842
$ws.Range("H83").Value = 4495
$ws.Range("J83").Value = 3990
$ws.Range("L83").Value = 35910
$ws.Range("N83").Value = -45270
$ws.Range("H113").Value = 1396.7142
$ws.Range("J113").Value = 746.5454999999999
$ws.Range("L113").Value = 2239.6365
$ws.Range("N113").Value = -6579.6365
$ws.Range("H122").Value = 1993
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1993
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 17937
$ws.Range("N122").Value = -22837
$ws.Range("M122").ClearContents()
$ws.Range("H131").Value = 6589167
$ws.Range("J131").Value = 10620.945
$ws.Range("L131").Value = 31862.835
$ws.Range("N131").Value = -41942.835
$ws.Range("H132").Value = 1596.8462
$ws.Range("J132").Value = 1729
$ws.Range("L132").Value = 15561
$ws.Range("N132").Value = -20621
$ws.Range("H134").Value = 2308.5557
$ws.Range("I134").Value = 2308.5557
$ws.Range("K134").Value = 6925.6671
$ws.Range("M134").Value = -1855.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 580.7
$ws.Range("I107").Value = 134
$ws.Range("J107").Value = 1250.75
$ws.Range("K107").Value = 134
$ws.Range("L107").Value = 1250.75
$ws.Range("M107").Value = 1786
$ws.Range("N107").Value = -5090.75
$ws.Range("H122").Value = 2482.5715
$ws.Range("I122").Value = 1681.6
$ws.Range("J122").Value = 2927.5557
$ws.Range("K122").Value = 5044.799999999999
$ws.Range("L122").Value = 8782.667099999999
$ws.Range("M122").Value = -2594.799999999999
$ws.Range("N122").Value = -13682.6671
$ws.Range("H126").Value = 2178060
$ws.Range("I126").Value = 2780727.2
$ws.Range("K126").Value = 8342181.600000001
$ws.Range("M126").Value = -8339711.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1620.2858
$ws.Range("I22").Value = 780.6667
$ws.Range("K22").Value = 780.6667
$ws.Range("M22").Value = -485.6667
$ws.Range("H27").Value = 1620.2858
$ws.Range("I27").Value = 780.6667
$ws.Range("K27").Value = 780.6667
$ws.Range("M27").Value = -673.6667
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H61").Value = 2361.3333
$ws.Range("I61").Value = 1892.8572
$ws.Range("J61").Value = 4001
$ws.Range("K61").Value = 1892.8572
$ws.Range("L61").Value = 4001
$ws.Range("M61").Value = -1690.8572
$ws.Range("N61").Value = -4405
$ws.Range("H93").Value = 1351.8
$ws.Range("I93").Value = 933
$ws.Range("K93").Value = 933
$ws.Range("M93").Value = 315
$ws.Range("H113").Value = 2361.3333
$ws.Range("I113").Value = 1892.8572
$ws.Range("J113").Value = 4001
$ws.Range("K113").Value = 1892.8572
$ws.Range("L113").Value = 4001
$ws.Range("M113").Value = 277.1428000000001
$ws.Range("N113").Value = -8341
$ws.Range("H122").Value = 6920.5
$ws.Range("I122").Value = 8252.857
$ws.Range("J122").Value = 6072.636
$ws.Range("K122").Value = 24758.571
$ws.Range("L122").Value = 18217.908
$ws.Range("M122").Value = -22308.571
$ws.Range("N122").Value = -23117.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2202.3076
$ws.Range("I81").Value = 1920.909
$ws.Range("K81").Value = 3841.818
$ws.Range("M81").Value = -2780.818
$ws.Range("H84").Value = 2202.3076
$ws.Range("I84").Value = 1920.909
$ws.Range("K84").Value = 19209.09
$ws.Range("M84").Value = -13905.09
$ws.Range("H122").Value = 107658.84
$ws.Range("J122").Value = 1947.8
$ws.Range("L122").Value = 5843.4
$ws.Range("N122").Value = -10743.4
